# Scheduled market-data refresh: update price/profit columns (H-N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 368168.34
$ws.Range("J13").Value = 2250
$ws.Range("L13").Value = 2250
$ws.Range("N13").Value = -2588
$ws.Range("H16").Value = 10000001
$ws.Range("I16").Value = 10000001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 10000001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = -9999771
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41944
$ws.Range("H54").Value = 15998.889
$ws.Range("I54").Value = 11198
$ws.Range("J54").Value = 22000
$ws.Range("K54").Value = 11198
$ws.Range("L54").Value = 22000
$ws.Range("M54").Value = -10712
$ws.Range("N54").Value = -22972
$ws.Range("H62").Value = 59827.5
$ws.Range("I62").Value = 94181.55
$ws.Range("K62").Value = 94181.55
$ws.Range("M62").Value = -93557.55
$ws.Range("H65").Value = 59827.5
$ws.Range("I65").Value = 94181.55
$ws.Range("K65").Value = 470907.75
$ws.Range("M65").Value = -467787.75
$ws.Range("H94").Value = 1542.8334
$ws.Range("I94").Value = 1542.8334
$ws.Range("K94").Value = 1542.8334
$ws.Range("M94").Value = -1091.8334
$ws.Range("H107").Value = 687.8125
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null
$ws.Range("H132").Value = 3488.7693
$ws.Range("I132").Value = 1218.9231
$ws.Range("J132").Value = 8028.4614
$ws.Range("K132").Value = 3656.7693
$ws.Range("L132").Value = 24085.3842
$ws.Range("M132").Value = -1126.7693
$ws.Range("N132").Value = -29145.3842
$ws.Range("H135").Value = 53535.74
$ws.Range("I135").Value = 894.25
$ws.Range("K135").Value = 8048.25
$ws.Range("M135").Value = -5513.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3100.158
$ws.Range("I2").Value = 2821.2666
$ws.Range("J2").Value = 4146
$ws.Range("K2").Value = 2821.2666
$ws.Range("L2").Value = 4146
$ws.Range("M2").Value = -2708.2666
$ws.Range("N2").Value = -4372
$ws.Range("H32").Value = 47081.2
$ws.Range("I32").Value = 28565.027
$ws.Range("K32").Value = 28565.027
$ws.Range("M32").Value = -28278.027
$ws.Range("H74").Value = 2237.3333
$ws.Range("I74").Value = 2237.3333
$ws.Range("K74").Value = 2237.3333
$ws.Range("M74").Value = -1363.3333
$ws.Range("H77").Value = 2237.3333
$ws.Range("I77").Value = 2237.3333
$ws.Range("K77").Value = 11186.6665
$ws.Range("M77").Value = -6818.666499999999
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984
$ws.Range("H116").Value = 3100.158
$ws.Range("I116").Value = 2821.2666
$ws.Range("J116").Value = 4146
$ws.Range("K116").Value = 2821.2666
$ws.Range("L116").Value = 4146
$ws.Range("M116").Value = -527.2665999999999
$ws.Range("N116").Value = -8734
$ws.Range("H132").Value = 2740.3704
$ws.Range("I132").Value = 1801.6154
$ws.Range("K132").Value = 5404.8462
$ws.Range("M132").Value = -2874.8462
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3100.158
$ws.Range("I3").Value = 2821.2666
$ws.Range("J3").Value = 4146
$ws.Range("K3").Value = 2821.2666
$ws.Range("L3").Value = 4146
$ws.Range("M3").Value = -2707.2666
$ws.Range("N3").Value = -4374
$ws.Range("H86").Value = 1909.0714
$ws.Range("I86").Value = 1588.6957
$ws.Range("K86").Value = 1588.6957
$ws.Range("M86").Value = -465.6957
$ws.Range("H89").Value = 1909.0714
$ws.Range("I89").Value = 1588.6957
$ws.Range("K89").Value = 7943.4785
$ws.Range("M89").Value = -2327.4785
$ws.Range("H107").Value = 14350.143
$ws.Range("I107").Value = 3915.5386
$ws.Range("K107").Value = 3915.5386
$ws.Range("M107").Value = -1995.5386
$ws.Range("H134").Value = 2731.913
$ws.Range("I134").Value = 2346.3333
$ws.Range("J134").Value = 3454.875
$ws.Range("K134").Value = 7038.999899999999
$ws.Range("L134").Value = 10364.625
$ws.Range("M134").Value = -4503.999899999999
$ws.Range("N134").Value = -15434.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3136.8
$ws.Range("I16").Value = 2561.3333
$ws.Range("K16").Value = 2561.3333
$ws.Range("M16").Value = -2274.3333
$ws.Range("H58").Value = 1738.9286
$ws.Range("I58").Value = 1531.3636
$ws.Range("K58").Value = 1531.3636
$ws.Range("M58").Value = -1328.3636
$ws.Range("H105").Value = 4189.3237
$ws.Range("I105").Value = 3822.4
$ws.Range("J105").Value = 4713.5
$ws.Range("K105").Value = 3822.4
$ws.Range("L105").Value = 4713.5
$ws.Range("M105").Value = -2075.4
$ws.Range("N105").Value = -8207.5
$ws.Range("H113").Value = 3136.8
$ws.Range("I113").Value = 2561.3333
$ws.Range("K113").Value = 2561.3333
$ws.Range("M113").Value = -391.3332999999998
$ws.Range("H132").Value = 2372.1
$ws.Range("I132").Value = 2391.9473
$ws.Range("K132").Value = 7175.841899999999
$ws.Range("M132").Value = -4645.841899999999
$ws.Range("H136").Value = 1738.9286
$ws.Range("I136").Value = 1531.3636
$ws.Range("K136").Value = 4594.0908
$ws.Range("M136").Value = -2044.0908
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 11500
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372
$ws.Range("H65").Value = 11500
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864
$ws.Range("H68").Value = 1346.2222
$ws.Range("I68").Value = 533.3333
$ws.Range("J68").Value = 1752.6666
$ws.Range("K68").Value = 1599.9999
$ws.Range("L68").Value = 5257.9998
$ws.Range("M68").Value = -788.9999
$ws.Range("N68").Value = -6879.9998
$ws.Range("H71").Value = 1346.2222
$ws.Range("I71").Value = 533.3333
$ws.Range("J71").Value = 1752.6666
$ws.Range("K71").Value = 4799.9997
$ws.Range("L71").Value = 15773.9994
$ws.Range("M71").Value = -743.9997000000003
$ws.Range("N71").Value = -23885.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1313.2142
$ws.Range("I102").Value = 838.5
$ws.Range("K102").Value = 838.5
$ws.Range("M102").Value = 783.5
$ws.Range("H132").Value = 1994.5
$ws.Range("I132").Value = 1965.625
$ws.Range("J132").Value = 2052.25
$ws.Range("K132").Value = 5896.875
$ws.Range("L132").Value = 6156.75
$ws.Range("M132").Value = -3366.875
$ws.Range("N132").Value = -11216.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 481.45456
$ws.Range("I16").Value = 481.45456
$ws.Range("K16").Value = 481.45456
$ws.Range("M16").Value = -311.45456
$ws.Range("H22").Value = 1549.8422
$ws.Range("I22").Value = 1545
$ws.Range("K22").Value = 1545
$ws.Range("M22").Value = -1250
$ws.Range("H27").Value = 1549.8422
$ws.Range("I27").Value = 1545
$ws.Range("K27").Value = 1545
$ws.Range("M27").Value = -1438
$ws.Range("H40").Value = 5560.8
$ws.Range("I40").Value = 5002.1665
$ws.Range("J40").Value = 6398.75
$ws.Range("K40").Value = 5002.1665
$ws.Range("L40").Value = 6398.75
$ws.Range("M40").Value = -4866.1665
$ws.Range("N40").Value = -6670.75
$ws.Range("H122").Value = 4634.6665
$ws.Range("I122").Value = 3952
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 11856
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -9406
$ws.Range("N122").Value = -22900
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = $null
$ws.Range("N134").Value = 0
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4249.75
$ws.Range("I136").Value = 2333
$ws.Range("K136").Value = 6999
$ws.Range("M136").Value = -4449
